$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values are purely numeric-looking strings (e.g. "234.15").
# Excel would normally auto-convert such text into a real number when assigned
# via .Value. The source data keeps these as plain text cells, so we temporarily
# force a text number format before assigning, then restore the default style.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.191.46'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '1.863.88'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '234.15'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("D7").Value = '0.4672'
$ws.Range("E7").Value = '  -0.71%  '

$ws.Range("D8").Value = '43.15'
$ws.Range("E8").Value = '  +0.54%  '

$ws.Range("D9").Value = '0.2853'
$ws.Range("E9").Value = '  -1.93%  '

$ws.Range("D10").Value = '0.06482'
$ws.Range("E10").Value = '  -2.22%  '

$ws.Range("D11").Value = '21.29'
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").Value = '0.07765'
$ws.Range("E12").Value = '  -3.13%  '

$ws.Range("D13").Value = '1.853.07'
$ws.Range("E13").Value = '  -0.89%  '

$ws.Range("D14").Value = '93.88'
$ws.Range("E14").Value = '  -3.66%  '

$ws.Range("D15").Value = '0.6826'
$ws.Range("E15").Value = '  -0.66%  '

$ws.Range("D16").Value = '5.050'
$ws.Range("E16").Value = '  -2.02%  '

$ws.Range("D17").Value = '271.46'
$ws.Range("E17").Value = '  -1.04%  '

$ws.Range("D18").Value = '30.168.51'
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").Value = '13.35'
$ws.Range("E19").Value = '  -5.23%  '

$ws.Range("D20").Value = '0.000007564'
$ws.Range("E20").Value = '  -2.19%  '

$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").Value = '  -0.17%  '

$ws.Range("D22").Value = '2.125.98'
$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("D23").Value = '0.9993'
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").Value = '5.145'
$ws.Range("E24").Value = '  -3.09%  '

$ws.Range("D25").Value = '6.112'
$ws.Range("E25").Value = '  -1.60%  '

$ws.Range("D26").Value = '9.365'
$ws.Range("E26").Value = '  +0.97%  '

$ws.Range("D27").Value = '165.48'
$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").Value = '18.56'
$ws.Range("E28").Value = '  -2.32%  '

$ws.Range("E29").Value = '  -3.57%  '

$ws.Range("D30").Value = '1.362'
$ws.Range("E30").Value = '  -0.49%  '

$ws.Range("D31").Value = '0.09839'
$ws.Range("E31").Value = '  -1.25%  '

$ws.Range("D32").Value = '1.453'
$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("D33").Value = '4.240'
$ws.Range("E33").Value = '  -2.85%  '

$ws.Range("D34").Value = '3.991'
$ws.Range("E34").Value = '  -2.32%  '

$ws.Range("E35").Value = '  -1.05%  '

$ws.Range("D36").Value = '1.118'
$ws.Range("E36").Value = '  -1.32%  '

$ws.Range("D37").Value = '0.6897'
$ws.Range("E37").Value = '  -2.04%  '

$ws.Range("D38").Value = '2.703'
$ws.Range("E38").Value = '  -0.20%  '

$ws.Range("D39").Value = '0.01836'
$ws.Range("E39").Value = '  -2.36%  '

$ws.Range("D40").Value = '2.752'
$ws.Range("E40").Value = '  +4.48%  '

$ws.Range("D41").Value = '6.362'
$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("E42").Value = '  -3.40%  '

$ws.Range("D43").Value = '0.9990'
$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("D44").Value = '1.885'
$ws.Range("E44").Value = '  -3.96%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8296'
$ws.Range("E45").Value = '  -1.45%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '102.42'
$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").Value = '0.4062'
$ws.Range("E47").Value = '  -2.56%  '

$ws.Range("D48").Value = '934.77'
$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("E49").Value = '  -1.95%  '

$ws.Range("D50").Value = '9.010'
$ws.Range("E50").Value = '  -3.09%  '

$ws.Range("D51").Value = '33.97'
$ws.Range("E51").Value = '  -1.45%  '

# Restore the default ("Normal") cell style on the forced-text cells so the
# workbook does not retain any extra text number-format styling.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
